$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": zero out species counts, drop the percentage column ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2:B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# --- Sheet "Species qualification": Range Analysis count drops to 0 ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# --- Sheet "High Priority break-up": re-summarised after mapping change ---
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Range("D2").Value = 3
$wsBreak.Range("E2").Value = 25

$wsBreak.Range("A3").Value = "IUCN"
$wsBreak.Range("B3").Value = 9
$wsBreak.Range("C3").Value = 75
$wsBreak.Range("D3").Value = 9
$wsBreak.Range("E3").Value = 75

$wsBreak.Rows.Item(4).Delete()
